# Fix the typo "data base" -> "database" in the project-objectives
# paragraph (keeps the rest of the sentence, including the other two
# correctly-spelled occurrences of "database" elsewhere in the doc,
# untouched by matching on the surrounding unique text).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "data base throught my sql",  # FindText
    $true,                        # MatchCase
    $false,                       # MatchWholeWord
    $false,                       # MatchWildcards
    $false,                       # MatchSoundsLike
    $false,                       # MatchAllWordForms
    $true,                        # Forward
    1,                            # Wrap (wdFindContinue)
    $false,                       # Format
    "database throught my sql",   # ReplaceWith
    2                             # Replace (wdReplaceAll)
)
